# May 9th update: new gyroscope readings (16 rows) pushed to the top of the
# x/y/z log, oldest readings pushed down; the bottom of the sliding window
# grows from 20 to 30 data rows (sheet dimension A1:C21 -> A1:C31).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.008163382383910248
$ws.Cells.Item(2, 2).Value = -0.07067988338795567
$ws.Cells.Item(2, 3).Value = -0.05231227912008755
$ws.Cells.Item(3, 1).Value = -0.06663984263485132
$ws.Cells.Item(3, 2).Value = -0.06482113017277281
$ws.Cells.Item(3, 3).Value = -0.02257424863902029
$ws.Cells.Item(4, 1).Value = -0.005289537662809526
$ws.Cells.Item(4, 2).Value = -0.05883742868900297
$ws.Cells.Item(4, 3).Value = 0.1704730540513992
$ws.Cells.Item(5, 1).Value = 0.00837163017554712
$ws.Cells.Item(5, 2).Value = -0.05375613881783047
$ws.Cells.Item(5, 3).Value = 0.2080968008799986
$ws.Cells.Item(6, 1).Value = -0.05617183616215521
$ws.Cells.Item(6, 2).Value = 0.01428591663187196
$ws.Cells.Item(6, 3).Value = 0.0969609718092462
$ws.Cells.Item(7, 1).Value = 0.0280998013913631
$ws.Cells.Item(7, 2).Value = 0.0762054398655891
$ws.Cells.Item(7, 3).Value = -0.0085521135479211
$ws.Cells.Item(8, 1).Value = 0.06225272437388239
$ws.Cells.Item(8, 2).Value = 0.09189357202161434
$ws.Cells.Item(8, 3).Value = -0.06575131636451573
$ws.Cells.Item(9, 1).Value = 0.04992435195229265
$ws.Cells.Item(9, 2).Value = 0.03760986301031981
$ws.Cells.Item(9, 3).Value = -0.06972193650223989
$ws.Cells.Item(10, 1).Value = 0.03176499361341648
$ws.Cells.Item(10, 2).Value = -0.02544809200546947
$ws.Cells.Item(10, 3).Value = 0.0004164989699015492
$ws.Cells.Item(11, 1).Value = 0.03686016493222929
$ws.Cells.Item(11, 2).Value = -0.01995030405338506
$ws.Cells.Item(11, 3).Value = 0.01384165091440073
$ws.Cells.Item(12, 1).Value = 0.01492454890500421
$ws.Cells.Item(12, 2).Value = 0.02186619833281094
$ws.Cells.Item(12, 3).Value = 0.002901610000249519
$ws.Cells.Item(13, 1).Value = -0.02151911604133509
$ws.Cells.Item(13, 2).Value = 0.04574547572569409
$ws.Cells.Item(13, 3).Value = -0.03826237769151855
$ws.Cells.Item(14, 1).Value = -0.01338350129398431
$ws.Cells.Item(14, 2).Value = 0.03186217543076381
$ws.Cells.Item(14, 3).Value = -0.08599316531961611
$ws.Cells.Item(15, 1).Value = 0.02753058617765249
$ws.Cells.Item(15, 2).Value = 0.01243943788788525
$ws.Cells.Item(15, 3).Value = -0.08503521572459828
$ws.Cells.Item(16, 1).Value = 0.02497605843977489
$ws.Cells.Item(16, 2).Value = 0.07183220271359811
$ws.Cells.Item(16, 3).Value = -0.08203642476688736
$ws.Cells.Item(17, 1).Value = 0.02700301970947864
$ws.Cells.Item(17, 2).Value = 0.2729595926674929
$ws.Cells.Item(17, 3).Value = -0.1162587641315027
$ws.Cells.Item(18, 1).Value = 0.0572686158120632
$ws.Cells.Item(18, 2).Value = 0.3712533414363861
$ws.Cells.Item(18, 3).Value = -0.058643065392971
$ws.Cells.Item(19, 1).Value = -0.1334879455918612
$ws.Cells.Item(19, 2).Value = 0.2170098776167092
$ws.Cells.Item(19, 3).Value = 0.09671107679605455
$ws.Cells.Item(20, 1).Value = -0.05985091050917481
$ws.Cells.Item(20, 2).Value = 0.1304891434582799
$ws.Cells.Item(20, 3).Value = -0.121909264813769
$ws.Cells.Item(21, 1).Value = -0.05091006410392845
$ws.Cells.Item(21, 2).Value = 0.1894376386295665
$ws.Cells.Item(21, 3).Value = -0.5621209686452691
$ws.Cells.Item(22, 1).Value = 0.0685279762203039
$ws.Cells.Item(22, 2).Value = -0.7451306256380917
$ws.Cells.Item(22, 3).Value = -1.788169232281768
$ws.Cells.Item(23, 1).Value = -0.745783193544906
$ws.Cells.Item(23, 2).Value = -2.762346484444355
$ws.Cells.Item(23, 3).Value = -1.756598580967298
$ws.Cells.Item(24, 1).Value = -0.1823294054378186
$ws.Cells.Item(24, 2).Value = -0.2232989397915999
$ws.Cells.Item(24, 3).Value = 1.138486255298958
$ws.Cells.Item(25, 1).Value = 0.5607048923319088
$ws.Cells.Item(25, 2).Value = 2.322676355188556
$ws.Cells.Item(25, 3).Value = 2.414278323000138
$ws.Cells.Item(26, 1).Value = -0.7038278525525871
$ws.Cells.Item(26, 2).Value = -0.9540881785479458
$ws.Cells.Item(26, 3).Value = -0.1982396515932949
$ws.Cells.Item(27, 1).Value = -0.2033903666517955
$ws.Cells.Item(27, 2).Value = -0.4635634205558086
$ws.Cells.Item(27, 3).Value = 0.1531466895883727
$ws.Cells.Item(28, 1).Value = 0.07737163657491758
$ws.Cells.Item(28, 2).Value = -0.3377668315714036
$ws.Cells.Item(28, 3).Value = 0.5267046527429067
$ws.Cells.Item(29, 1).Value = 0.437226802110672
$ws.Cells.Item(29, 2).Value = -1.076497316360474
$ws.Cells.Item(29, 3).Value = 0.305127203464508
$ws.Cells.Item(30, 1).Value = 0.3028364533727822
$ws.Cells.Item(30, 2).Value = -0.5765596194700771
$ws.Cells.Item(30, 3).Value = -0.2633939818902439
$ws.Cells.Item(31, 1).Value = -0.01123159040104305
$ws.Cells.Item(31, 2).Value = -0.1864527369087388
$ws.Cells.Item(31, 3).Value = -0.3065016323869878
